$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string reorderings (country name swaps) ---
# Row 80/81 swap: "Estado de Palestina" <-> "Bulgaria" (Bulgaria moves up with updated stats;
# Palestina keeps its former row-80 numbers, now on row 81).
$ws.Range("A80").Value = "Bulgaria"
$ws.Range("A81").Value = "Estado de Palestina"

# Row 210/211 swap: "Groenlandia" <-> "Islas Malvinas" (identical stats, so only labels move).
$ws.Range("A210").Value = "Islas Malvinas"
$ws.Range("A211").Value = "Groenlandia"

# --- Updated case counts for the new report timestamp ---
# Row 4
$ws.Range("B4").Value = 4237908
$ws.Range("C4").Value = 67590
$ws.Range("D4").Value = 2010465
$ws.Range("E4").Value = 2079162
$ws.Range("G4").Value = 932
$ws.Range("H4").Value = 148281

# Row 5
$ws.Range("B5").Value = 2343366
$ws.Range("C5").Value = 53415
$ws.Range("D5").Value = 1592281
$ws.Range("E5").Value = 665847
$ws.Range("G5").Value = 1031
$ws.Range("H5").Value = 85238

# Row 9
$ws.Range("B9").Value = 375961
$ws.Range("C9").Value = 4865
$ws.Range("D9").Value = 259423
$ws.Range("E9").Value = 98695
$ws.Range("G9").Value = 189
$ws.Range("H9").Value = 17843

# Row 21
$ws.Range("B21").Value = 205960
$ws.Range("C21").Value = 818
$ws.Range("E21").Value = 7359
$ws.Range("G21").Value = 14
$ws.Range("H21").Value = 9201

# Row 51
$ws.Range("B51").Value = 38458
$ws.Range("C51").Value = 462
$ws.Range("D51").Value = 34826
$ws.Range("E51").Value = 3496

# Row 80
$ws.Range("B80").Value = 10123
$ws.Range("C80").Value = 270
$ws.Range("D80").Value = 5252
$ws.Range("E80").Value = 4534
$ws.Range("G80").Value = 8
$ws.Range("H80").Value = 337

# Row 81
$ws.Range("B81").Value = 10093
$ws.Range("C81").Value = 349
$ws.Range("D81").Value = 2720
$ws.Range("E81").Value = 7303
$ws.Range("G81").Value = 3
$ws.Range("H81").Value = 70

# Row 130
$ws.Range("B130").Value = 1765
$ws.Range("C130").Value = 13
$ws.Range("D130").Value = 1297
$ws.Range("E130").Value = 402

# Row 186
$ws.Range("B186").Value = 108
$ws.Range("C186").Value = 2
$ws.Range("E186").Value = 7

# Row 214
$ws.Range("B214").Value = 11
$ws.Range("C214").Value = 1
$ws.Range("E214").Value = 4

# --- Report timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 25 de Julio de 2020 a las 00:01"
